$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 8 (home_xGoals), shifting
# everything from the old row 8 down to row 12.
$ws.Rows.Item(8).Resize(4).Insert()

# Copy the label-column formatting (bold, centered, bordered style used by
# column A headers) down into the newly inserted label cells.
$ws.Range("A12").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)

# Populate the four new rows with the goals-related stats.
$ws.Range("A8").Value = "homeGoals"
$ws.Range("B8").Value = 10

$ws.Range("A9").Value = "awayGoals"
$ws.Range("B9").Value = 9

$ws.Range("A10").Value = "homeGoalsHalfTime"
$ws.Range("B10").Value = 6

$ws.Range("A11").Value = "awayGoalsHalfTime"
$ws.Range("B11").Value = 5
